$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as literal text (never auto-converted to a
# number/date by Excel's COM value-coercion), matching the inlineStr cells in the
# source workbook, and without leaving a stray "Text" cell style behind.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "47.262.68"
$ws.Range("E2").Value = "  +4.59%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.483.55"
$ws.Range("E3").Value = "  +2.29%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
Set-TextValue $ws.Range("D5") "322.74"
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
Set-TextValue $ws.Range("D6") "104.60"
$ws.Range("E6").Value = "  +1.40%  "

# Row 7
$ws.Range("E7").Value = "  +1.30%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  -0.13%  "

# Row 9
$ws.Range("E9").Value = "  +2.12%  "

# Row 10
Set-TextValue $ws.Range("D10") "37.28"
$ws.Range("E10").Value = "  +4.55%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0811"
$ws.Range("E11").Value = "  +1.14%  "

# Row 12
$ws.Range("E12").Value = "  +0.14%  "

# Row 13
Set-TextValue $ws.Range("D13") "18.27"
$ws.Range("E13").Value = "  -0.71%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.16"
$ws.Range("E14").Value = "  +2.66%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.872.40"
$ws.Range("E15").Value = "  +2.11%  "

# Row 16
Set-TextValue $ws.Range("D16") "2.479.84"
$ws.Range("E16").Value = "  +1.82%  "

# Row 17
$ws.Range("E17").Value = "  +1.35%  "

# Row 18
Set-TextValue $ws.Range("D18") "47.148.86"
$ws.Range("E18").Value = "  +4.60%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.67"
$ws.Range("E19").Value = "  +3.43%  "

# Row 20
$ws.Range("E20").Value = "  +2.69%  "

# Row 21
$ws.Range("E21").Value = "  +1.11%  "

# Row 22
Set-TextValue $ws.Range("D22") "70.55"
$ws.Range("E22").Value = "  +2.35%  "

# Row 23
Set-TextValue $ws.Range("D23") "250.17"
$ws.Range("E23").Value = "  +2.77%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.36"
$ws.Range("E24").Value = "  +4.46%  "

# Row 25
$ws.Range("E25").Value = "  +2.17%  "

# Row 26
Set-TextValue $ws.Range("D26") "26.09"
$ws.Range("E26").Value = "  +3.06%  "

# Row 27
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D28") "10.07"
$ws.Range("E28").Value = "  +5.71%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D29") "2.20"
$ws.Range("E29").Value = "  -2.82%  "

# Row 30
Set-TextValue $ws.Range("D30") "35.18"
$ws.Range("E30").Value = "  +6.75%  "

# Row 32
Set-TextValue $ws.Range("D32") "49.53"
$ws.Range("E32").Value = "  +0.78%  "

# Row 33
$ws.Range("E33").Value = "  -1.02%  "

# Row 34
$ws.Range("E34").Value = "  +2.91%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.0780"
$ws.Range("E35").Value = "  +2.08%  "

# Row 36
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
Set-TextValue $ws.Range("D37") "4.64"
$ws.Range("E37").Value = "  +3.63%  "

# Row 38
$ws.Range("E38").Value = "  +2.35%  "

# Row 39
$ws.Range("E39").Value = "  +4.49%  "

# Row 40
$ws.Range("E40").Value = "  +1.55%  "

# Row 41
Set-TextValue $ws.Range("D41") "121.17"
$ws.Range("E41").Value = "  -2.46%  "

# Row 42
$ws.Range("E42").Value = "  +0.90%  "

# Row 43
Set-TextValue $ws.Range("D43") "21.38"
$ws.Range("E43").Value = "  +0.45%  "

# Row 44
$ws.Range("E44").Value = "  +1.78%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.947.90"
$ws.Range("E45").Value = "  +0.64%  "

# Row 46
$ws.Range("E46").Value = "  +1.64%  "

# Row 47
$ws.Range("E47").Value = "  +0.05%  "

# Row 48
$ws.Range("E48").Value = "  -0.62%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.79"
$ws.Range("E49").Value = "  +0.85%  "

# Row 50
Set-TextValue $ws.Range("D50") "5.36"
$ws.Range("E50").Value = "  +13.42%  "

# Row 51
Set-TextValue $ws.Range("D51") "78.53"
$ws.Range("E51").Value = "  +3.12%  "

